$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Formula = "=""2021-07-13"""
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial(-4163)

$ws.Range("D14").Formula = "=""16:07:33"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)

$ws.Range("C15").Formula = "=""2021-07-13"""
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

$ws.Range("D15").Formula = "=""16:08:18"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
